# Auto update Excel log
# Appends newly-logged sensor readings to the PIR, Humidity and Temperature
# sheets (mirrors the source system's periodic CSV->Excel export).

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [object[]]$Rows
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $endRow = $StartRow + $Rows.Count - 1

    # Columns A (Date) and E (Value, e.g. "88.4%") look numeric to Excel's
    # smart-typing on assignment; force them to Text first so the values are
    # stored verbatim as strings, matching the rest of the log.
    $ws.Range("A$StartRow`:A$endRow").NumberFormat = "@"
    $ws.Range("E$StartRow`:E$endRow").NumberFormat = "@"

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $rowData = $Rows[$i]
        for ($j = 0; $j -lt $rowData.Count; $j++) {
            $c = $j + 1
            $ws.Cells.Item($r, $c).Value = $rowData[$j]
        }
    }
}

# PIR sheet: 13 new "No Motion" / "Inactive" readings (rows 14-26)
$pirRows = @(
    @("2026-01-28", "16:11:32", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:11:33", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:11:38", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:11:43", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:11:48", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:11:53", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:11:58", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:12:03", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:12:08", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:12:13", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:12:18", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:12:23", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:12:28", "16:00", "Bathroom", "No Motion", "Inactive")
)
Append-Rows "PIR" 14 $pirRows

# Humidity sheet: 11 new readings (rows 15-25)
$humidityRows = @(
    @("2026-01-28", "16:11:33", "16:00", "Bathroom", "88.4%", "Active"),
    @("2026-01-28", "16:11:37", "16:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "16:11:45", "16:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "16:11:49", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:11:53", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:12:01", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:12:09", "16:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "16:12:13", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:12:17", "16:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "16:12:21", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:12:29", "16:00", "Bathroom", "87.4%", "Active")
)
Append-Rows "Humidity" 15 $humidityRows

# Temperature sheet: 11 new readings (rows 15-25)
$temperatureRows = @(
    @("2026-01-28", "16:11:33", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:11:37", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:11:45", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:11:49", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:11:53", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:12:01", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:12:09", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:12:13", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:12:17", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:12:21", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:12:29", "16:00", "Bathroom", "22.8C", "Active")
)
Append-Rows "Temperature" 15 $temperatureRows
